# Feedback tracker: mark the addressed items and flag the two that still
# need a follow-up discussion ("Rücksprache") before they can be closed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8 already carries the "addressed" (green) fill style from before - reuse
# it as the template for the other rows whose feedback has now been handled.
$ws.Cells.Item(8, 2).Copy()

$addressedRows = @(2, 4, 12, 14, 16, 18, 24, 26, 28, 30, 32)
foreach ($r in $addressedRows) {
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)  # xlPasteFormats
}

$excel.CutCopyMode = 0

# These two items need further discussion before they can be marked done.
$ws.Cells.Item(6, 2).Value = "Rücksprache"
$ws.Cells.Item(10, 2).Value = "Rücksprache"

# Restore the selection/cursor position as left by the author.
$ws.Range("E8").Select()
